# The deck's applied design ("Integral" / "Red Violet" colour scheme) is
# switched back to the standard Office colour scheme, i.e. each of the
# twelve theme colour slots (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) on the Slide Master's theme is re-pointed at the default
# Office palette.

$p = $ppt.ActivePresentation

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target ("Office") palette, in ThemeColorScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    (RGBVal 0x00 0x00 0x00),  # dk1      000000
    (RGBVal 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (RGBVal 0x44 0x54 0x6A),  # dk2      44546A
    (RGBVal 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (RGBVal 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (RGBVal 0xED 0x7D 0x31),  # accent2  ED7D31
    (RGBVal 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (RGBVal 0xFF 0xC0 0x00),  # accent4  FFC000
    (RGBVal 0x44 0x72 0xC4),  # accent5  4472C4
    (RGBVal 0x70 0xAD 0x47),  # accent6  70AD47
    (RGBVal 0x05 0x63 0xC1),  # hlink    0563C1
    (RGBVal 0x95 0x4F 0x72)   # folHlink 954F72
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
